# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header suffixes to the actual format-version
# identifiers ("_FV2310" / "_FV2404"), wraps the data range in a real Excel
# Table (ListObject) with an AutoFilter, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) ---------------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2310"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2404"
        }
    }
}

# --- 2. Turn the data range into a real table (ListObject) ----------------
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add(1, $headerRange, 0, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
